# Max News 4 -> 10
# 기사 기존 4개에서 최대 10개까지 출력 가능
#
# 1) Un-hide the "paper" and "internet" report sheets.
# 2) Extend the formula-driven "paper"/"internet" display sheets from 4
#    rows to 10 rows (rows 5-10 now pull from papPaste/intPaste just like
#    rows 1-4 already did).
# 3) Extend the "intPaste" paste-staging sheet from 4 news slots (rows
#    1-8) to 8 news slots (rows 1-16), duplicating the last slot's
#    layout/merges for the new rows.

$wb = $excel.ActiveWorkbook

# --- 1) Un-hide paper / internet ------------------------------------------------
$wb.Worksheets("paper").Visible = -1
$wb.Worksheets("internet").Visible = -1

# --- 2) paper sheet: fill rows 5-10 with the same OFFSET formulas as 1-4 -------
$paper = $wb.Worksheets("paper")
for ($row = 5; $row -le 10; $row++) {
    $paper.Range("A$row").Formula = "=IF(OFFSET(papPaste!`$A`$1,ROW(A$row)*2-2,0)=0,`"`",OFFSET(papPaste!`$A`$1,ROW(A$row)*2-2,0))"
    $paper.Range("B$row").Formula = "=IF(OFFSET(papPaste!`$A`$1,ROW(A$row)*2-2,1)=0,`"`",OFFSET(papPaste!`$A`$1,ROW(A$row)*2-2,1))"
    $paper.Range("C$row").Formula = "=IF(OFFSET(papPaste!`$A`$1,ROW(B$row)*2-2,2)=0,`"`",OFFSET(papPaste!`$A`$1,ROW(B$row)*2-2,2))"
    $paper.Range("D$row").Formula = "=IF(OFFSET(papPaste!`$A`$1,ROW(B$row)*2-2,3)=0,`"`",OFFSET(papPaste!`$A`$1,ROW(B$row)*2-2,3))"
    $paper.Range("E$row").Formula = "=IF(OFFSET(papPaste!`$A`$1,ROW(B$row)*2-1,2)=0,`"`",OFFSET(papPaste!`$A`$1,ROW(B$row)*2-1,2))"
}
$paper.Range("E5").Select()

# --- 2b) internet sheet: fill rows 5-10 with the same OFFSET formulas as 1-4 --
$internet = $wb.Worksheets("internet")
for ($row = 5; $row -le 10; $row++) {
    $internet.Range("A$row").Formula = "=IF(OFFSET(intPaste!`$A`$1,ROW(A$row)*2-2,0)=0,`"`",OFFSET(intPaste!`$A`$1,ROW(A$row)*2-2,0))"
    $internet.Range("B$row").Formula = "=IF(OFFSET(intPaste!`$A`$1,ROW(A$row)*2-2,1)=0,`"`",OFFSET(intPaste!`$A`$1,ROW(A$row)*2-2,1))"
    $internet.Range("C$row").Formula = "=IF(OFFSET(intPaste!`$A`$1,ROW(B$row)*2-2,2)=0,`"`",OFFSET(intPaste!`$A`$1,ROW(B$row)*2-2,2))"
    $internet.Range("D$row").Formula = "=IF(OFFSET(intPaste!`$A`$1,ROW(B$row)*2-2,3)=0,`"`",OFFSET(intPaste!`$A`$1,ROW(B$row)*2-2,3))"
    $internet.Range("E$row").Formula = "=IF(OFFSET(intPaste!`$A`$1,ROW(B$row)*2-1,2)=0,`"`",OFFSET(intPaste!`$A`$1,ROW(B$row)*2-1,2))"
}
$internet.Range("E5:E10").Select()

# --- 3) intPaste sheet: add 4 more two-row news slots (rows 9-16) -------------
# Each slot duplicates the 4th slot's (rows 7-8) content as a placeholder,
# keeping the same centered-alignment style on columns A/B/D.
$intPaste = $wb.Worksheets("intPaste")
$slotDate = "2022.09.30."
$slotLabel = "뉴스4"
$slotTitle = "뉴스4제목"
$slotBody = "뉴스4내용"
$slotLink = "뉴스4링크"

for ($slot = 0; $slot -lt 4; $slot++) {
    $topRow = 9 + ($slot * 2)
    $botRow = $topRow + 1

    $intPaste.Range("A$topRow").Value = $slotDate
    $intPaste.Range("A$topRow").HorizontalAlignment = -4108
    $intPaste.Range("B$topRow").Value = $slotLabel
    $intPaste.Range("B$topRow").HorizontalAlignment = -4108
    $intPaste.Range("C$topRow").Value = $slotTitle
    $intPaste.Range("D$topRow").Value = $slotBody
    $intPaste.Range("D$topRow").HorizontalAlignment = -4108

    $intPaste.Range("A$botRow").HorizontalAlignment = -4108
    $intPaste.Range("B$botRow").HorizontalAlignment = -4108
    $intPaste.Range("C$botRow").Value = $slotLink
    $intPaste.Range("D$botRow").HorizontalAlignment = -4108

    $intPaste.Range("A$topRow`:A$botRow").Merge()
    $intPaste.Range("B$topRow`:B$botRow").Merge()
    $intPaste.Range("D$topRow`:D$botRow").Merge()
}

$intPaste.Range("F10").Select()
